$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$range = $ws.Range("A1:U79")
$tbl = $ws.ListObjects.Add()
Write-Host "done"
